$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.732.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.87%  '

$ws.Range("D3").Value = "'2.364.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.76%  '

$ws.Range("E4").Value = '  -0.23%  '

$ws.Range("D5").Value = "'330.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.49%  '

$ws.Range("D6").Value = "'100.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.63%  '

$ws.Range("D7").Value = "'0.638"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.83%  '

$ws.Range("D9").Value = "'0.634"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.67%  '

$ws.Range("D10").Value = "'39.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.74%  '

$ws.Range("D11").Value = "'0.0922"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.60%  '

$ws.Range("D12").Value = "'8.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.70%  '

$ws.Range("E13").Value = '  -3.32%  '

$ws.Range("E14").Value = '  +0.28%  '

$ws.Range("D15").Value = "'16.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.98%  '

$ws.Range("D16").Value = "'2.719.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.67%  '

$ws.Range("D17").Value = "'2.362.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.74%  '

$ws.Range("D18").Value = "'42.666.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.97%  '

$ws.Range("D19").Value = "'7.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.28%  '

$ws.Range("D20").Value = "'0.0000107"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.90%  '

$ws.Range("D21").Value = "'3.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +9.85%  '

$ws.Range("D22").Value = "'75.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.60%  '

$ws.Range("D23").Value = "'269.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.75%  '

$ws.Range("E24").Value = '  -10.18%  '

$ws.Range("D25").Value = "'9.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.29%  '

$ws.Range("D26").Value = "'1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.58%  '

$ws.Range("D27").Value = "'11.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.48%  '

$ws.Range("D28").Value = "'23.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.31%  '

$ws.Range("D29").Value = "'2.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.02%  '

$ws.Range("D30").Value = "'176.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.19%  '

$ws.Range("D31").Value = "'3.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.87%  '

$ws.Range("D32").Value = "'0.0905"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.03%  '

$ws.Range("D33").Value = "'35.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -9.44%  '

$ws.Range("D34").Value = "'6.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.05%  '

$ws.Range("E35").Value = '  -0.49%  '

$ws.Range("D36").Value = "'4.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.74%  '

$ws.Range("E37").Value = '  -4.72%  '

$ws.Range("D38").Value = "'2.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.74%  '

$ws.Range("E39").Value = '  +2.01%  '

$ws.Range("D40").Value = "'3.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.49%  '

$ws.Range("D41").Value = "'1.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.21%  '

$ws.Range("D42").Value = "'0.235"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.18%  '

$ws.Range("D43").Value = "'70.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.42%  '

$ws.Range("E44").Value = '  -0.09%  '

$ws.Range("D45").Value = "'118.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.79%  '

$ws.Range("D46").Value = "'90.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +27.99%  '

$ws.Range("D47").Value = "'11.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.50%  '

$ws.Range("D48").Value = "'5.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.51%  '

$ws.Range("D49").Value = "'9.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.53%  '

$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = "'1.572.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.47%  '

$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").Value = "'1.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.30%  '
